$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1.xml)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 57
$ws.Range("D3").Value = 70
$ws.Range("F3").Value = 64
$ws.Range("G3").Value = 56
$ws.Range("J3").Value = 88
$ws.Range("K3").Value = 102
$ws.Range("C6").Value = 225
$ws.Range("F6").Value = 244
$ws.Range("G6").Value = 228
$ws.Range("H6").Value = 199
$ws.Range("I6").Value = 267
$ws.Range("K6").Value = 230
$ws.Range("C7").Value = 304
$ws.Range("D7").Value = 320
$ws.Range("F7").Value = 347
$ws.Range("G7").Value = 336
$ws.Range("H7").Value = 311
$ws.Range("I7").Value = 416
$ws.Range("J7").Value = 354
$ws.Range("K7").Value = 417

# By Neighborhood (sheet2.xml)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G4").Value = 3
$ws.Range("G9").Value = 3
$ws.Range("F25").Value = 6
$ws.Range("H25").Value = 7
$ws.Range("D26").Value = 22
$ws.Range("H26").Value = 28
$ws.Range("J26").Value = 13
$ws.Range("K34").Value = 31
$ws.Range("K43").Value = 4
$ws.Range("F49").Value = 4
$ws.Range("K51").Value = 57
$ws.Range("C61").Value = 2
$ws.Range("G67").Value = 2
$ws.Range("I68").Value = 8
$ws.Range("G74").Value = 9
$ws.Range("C95").Value = 304
$ws.Range("D95").Value = 320
$ws.Range("F95").Value = 347
$ws.Range("G95").Value = 336
$ws.Range("H95").Value = 311
$ws.Range("I95").Value = 416
$ws.Range("J95").Value = 354
$ws.Range("K95").Value = 417

# Roseland (sheet4.xml)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 9

# Grand Crossing (sheet11.xml)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 6
$ws.Range("K5").Value = 19
$ws.Range("K6").Value = 31

# Armour Square (sheet12.xml)
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 3

# Englewood (sheet18.xml)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H2").Value = 9
$ws.Range("D3").Value = 10
$ws.Range("J3").Value = 3
$ws.Range("D6").Value = 22
$ws.Range("H6").Value = 28
$ws.Range("J6").Value = 13

# Loop (sheet22.xml)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 57

# Avondale (sheet42.xml)
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3

# New City (sheet43.xml)
$ws = $wb.Worksheets.Item('New City')
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2

# Edgewater (sheet44.xml)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = 7
$ws.Range("F5").Value = 6
$ws.Range("H5").Value = 7

# Little Village (sheet49.xml)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("D3").Value = 1
$ws.Range("D6").Value = 4

# Old Town (sheet70.xml)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("H4").Value = 7
$ws.Range("H5").Value = 8

# Jefferson Park (sheet77.xml)
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 4

# Oakland (sheet81.xml)
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("D3").Value = 2
$ws.Range("D5").Value = 2
